$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A.
# This shifts: A->B (segment names), B->C (RawActivations),
#              C->D (PercActivations), D->E (totalActivation)
$ws.Columns.Item(1).Insert()

# New header cell for the inserted column
$ws.Range("B1").Value = "segments"

# Match the header formatting (bold / border / centered) used by the
# other header cells, by copying format from the neighbouring header cell.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# The new column A will hold the numeric segment index (0-based) and should
# carry the same style the label column used to have, so copy that format
# over before clearing it from the (now relocated) label column.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2:A20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# The label column (now B) is plain text with no special formatting.
$ws.Range("B2:B20").ClearFormats()

# Fill in the new segment-index column with 0-based numeric ids.
for ($i = 0; $i -le 18; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}
